## V 0.53-B50 - final testing
## Add ADF(1 only) Items Freq, Name, ID + Add Distance to Destination item
##
## Tabelle2 is the "EngineMerge" mapping sheet: row 1 holds the field-name
## headers, rows 2-40 hold per-aircraft-profile placeholder/override data,
## and the last two columns are fixed marker columns (END_OF_COL, Title).
## Four brand new fields (GPS_DST, ADF1, ADF1_NAME, ADF1_F) are inserted
## right before those trailing marker columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle2")
$ws.Activate()

# The sheet used to end with ... COMPASS | END_OF_COL | Title (cols EE:EG).
# Insert 4 fresh columns ahead of the END_OF_COL marker (EF:EI) - this
# shifts END_OF_COL/Title from EF/EG out to EJ/EK and slides every
# formula/reference along with them automatically.
$ws.Columns("EF:EI").Insert()

# New header row-1 labels for the inserted columns.
$ws.Range("EF1").Value = "GPS_DST"
$ws.Range("EG1").Value = "ADF1"
$ws.Range("EH1").Value = "ADF1_NAME"
$ws.Range("EI1").Value = "ADF1_F"

# Data rows 2-40 get the same "|" placeholder used by every other
# not-yet-configured feature column on this template sheet.
$ws.Range("EF2:EI40").Value = "|"

# Leave the cursor where the author left it after finishing the edit.
[void]$ws.Range("EL34").Select()
